# Generate Report for Archive
#
# The "Status" value shown for the 84ac6ebe-... file changes from
# "Ready for handoff" to "In Translation" everywhere it is reported:
#   - Overview sheet: columns "zh-cn" (E2) and "de-de" (F2)
#   - zh-cn sheet: "Status" column (C2)
#   - de-de sheet: "Status" column (C2)
#
# Because the new status text is shorter than the old one, the
# "Status"/"zh-cn"/"de-de" columns that were auto-sized to fit it become
# narrower as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the reported status text.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the columns that hold the status text to their new (narrower)
# content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
